# Atualização de bases das ligas, do dia: 17-06-2024 às 21:10
#
# The underlying source data for several fixtures got their rows mixed up
# (B:AD swapped between two/three rows that share the same match date).
# This script restores the correct pairing by moving the B:AD payload
# (id .. PL_AhUnder) between the affected rows while leaving column A
# (the row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowPayload($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# Simple 2-row swaps (same match date, rows got transposed)
Swap-RowPayload 20 21
Swap-RowPayload 65 66
Swap-RowPayload 119 120
Swap-RowPayload 215 216

# 3-row rotation: new(95) = old(97), new(96) = old(95), new(97) = old(96)
$r95 = $ws.Range("B95:AD95")
$r96 = $ws.Range("B96:AD96")
$r97 = $ws.Range("B97:AD97")

$v95 = $r95.Value2
$v96 = $r96.Value2
$v97 = $r97.Value2

$r95.Value2 = $v97
$r96.Value2 = $v95
$r97.Value2 = $v96
